$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("D29").Value = 44483
$ws.Range("H29").Value = "Dulce o Americano"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 450
$ws.Range("K29").Value = 25000
$ws.Range("L29").Value = 26000
$ws.Range("M29").Value = 25556
$ws.Range("N29").Value = "`$/malla 70 unidades"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 365
$ws.Range("Q29").Value = 70

# Row 30
$ws.Range("D30").Value = 44509
$ws.Range("H30").Value = "Dulce o Americano"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 20000
$ws.Range("L30").Value = 22000
$ws.Range("M30").Value = 21000
$ws.Range("N30").Value = "`$/malla 70 unidades"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 300
$ws.Range("Q30").Value = 70

# Row 31
$ws.Range("D31").Value = 44230
$ws.Range("H31").Value = "Choclero"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 320
$ws.Range("L31").Value = 350
$ws.Range("M31").Value = 335
$ws.Range("N31").Value = "`$/unidad"
$ws.Range("O31").Value = "Región de O'Higgins"
$ws.Range("P31").Value = 335
$ws.Range("Q31").Value = 1

# Row 32
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 280
$ws.Range("L32").Value = 280
$ws.Range("M32").Value = 280
$ws.Range("P32").Value = 280

# Row 33
$ws.Range("D33").Value = 44265
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 15000
$ws.Range("K33").Value = 200
$ws.Range("L33").Value = 250
$ws.Range("M33").Value = 233
$ws.Range("P33").Value = 233

# Row 34
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 150
$ws.Range("L34").Value = 150
$ws.Range("M34").Value = 150
$ws.Range("P34").Value = 150

# Row 35
$ws.Range("D35").Value = 44349
$ws.Range("H35").Value = "Dulce o Americano"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9600
$ws.Range("N35").Value = "`$/malla 60 unidades"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 160
$ws.Range("Q35").Value = 60

# Row 36
$ws.Range("D36").Value = 44398
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 24000
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = 24500
$ws.Range("P36").Value = 408

# Row 37
$ws.Range("D37").Value = 44232
$ws.Range("H37").Value = "Choclero"
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 320
$ws.Range("L37").Value = 350
$ws.Range("M37").Value = 335
$ws.Range("N37").Value = "`$/unidad"
$ws.Range("O37").Value = "Región de O'Higgins"
$ws.Range("P37").Value = 335
$ws.Range("Q37").Value = 1

# Row 38
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 280
$ws.Range("L38").Value = 280
$ws.Range("M38").Value = 280
$ws.Range("P38").Value = 280

# Row 39
$ws.Range("D39").Value = 44586
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 150
$ws.Range("L39").Value = 200
$ws.Range("M39").Value = 175
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 175

# Row 40
$ws.Range("D40").Value = 44589
$ws.Range("O40").Value = "Provincia de Chacabuco"
